# Student "Riyas ahamed J" (row 26) finished their last pending task
# ("create_canva-menu"), so:
#   - the Pending Task cell (D26) is cleared
#   - the Completion Status cell (E26) flips from "Pending" to "Completed"
#   - the row's Name/Status cells pick up the green "Completed" styling
#     already used elsewhere in the sheet (bold white font on a green fill)
# The report footer's "Generated: ..." timestamp is also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing "Completed" row formatting (bold white on green fill)
# onto this row's Name/Status cells, matching B2/E2 etc., instead of
# re-building the formatting property by property (which would mint new,
# unused style entries).
$ws.Range("B2").Copy()
$ws.Range("B26").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E26").PasteSpecial(-4122)

# Clear the now-finished pending task.
$ws.Range("D26").ClearContents()

# Flip the completion status for this student.
$ws.Range("E26").Value = "Completed"

# Refresh the "Generated:" timestamp in the footer.
$ws.Range("A29").Value = "Generated: 2023-09-02 12:00:19 PM"
